# Update daily power records
# - Row 56: correct the End Time (C56)
# - Row 57: fill in the End Time (C57), which was previously left blank
# - Row 58: add a new daily record, letting the table's calculated
#   columns (Duration / Second Duration / Absolute Value) populate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the End Time for 2018-10-07 (row 56)
$ws.Range("C56").Value = 0.875

# 2. Fill in the End Time for 2018-10-08 (row 57), previously blank
$ws.Range("C57").Value = 0

# 3. Grow the table by one row (A1:F57 -> A1:F58) for the new day's record
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# 4. Enter the new day's date (2018-10-09) and calculated-column formulas
$ws.Range("A58").Value = 43382
$ws.Range("D58").Formula = "=(C58-B58)* 1440"
$ws.Range("E58").Formula = "=IF(C58>B58, (C58-B58)*1440, (B58-C58)*1440)"
$ws.Range("F58").Formula = "=ABS((C58-B58)*1440)"
$ws.Range("E58").Style = $ws.Range("E57").Style
$ws.Range("F58").Style = $ws.Range("F57").Style

# 5. Update the view: scroll down a row and select the next entry cell
$win = $excel.ActiveWindow
$win.ScrollRow = 47
$win.ScrollColumn = 1
$ws.Range("B58").Select()
